$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 24: Sub_component (E24) changes from "drop" to "Information_Schema" ----
$ws.Range("E24").Value = "Information_Schema"

# ---- New rows 25-32 ----

# Row 25
$ws.Range("A25").Value = "ddl_024"
$ws.Range("B25").Value = "y"
$ws.Range("C25").Value = "创建database,验证元数据信息保存"
$ws.Range("D25").Value = "DDL"
$ws.Range("E25").Value = "databaseCreate"
$ws.Range("H25").Value = "create database MYDDL_024"
$ws.Range("I25").Value = "select * from information_schema.schemata where schema_name in ('META','DINGO','ROOT','MYSQL','INFORMATION_SCHEMA','MYDDL_024')"
$ws.Range("J25").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_024.csv"
$ws.Range("K25").Value = "csv_containsAll"

# Row 26
$ws.Range("A26").Value = "ddl_025"
$ws.Range("B26").Value = "y"
$ws.Range("C26").Value = "创建schema,验证元数据信息保存"
$ws.Range("D26").Value = "DDL"
$ws.Range("E26").Value = "schemaCreate"
$ws.Range("H26").Value = "create schema MYDDL_025"
$ws.Range("I26").Value = "select * from information_schema.schemata where schema_name in ('META','DINGO','ROOT','MYSQL','INFORMATION_SCHEMA','MYDDL_025')"
$ws.Range("J26").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_025.csv"
$ws.Range("K26").Value = "csv_containsAll"

# Row 27 (tall row, wrapped text)
$ws.Range("A27").Value = "ddl_026"
$ws.Range("B27").Value = "y"
$ws.Range("C27").Value = "删除database,验证元数据信息删除"
$ws.Range("D27").Value = "DDL"
$ws.Range("E27").Value = "databaseDrop"
$ws.Range("H27").Value = "create database MYDDL_026;drop database MYDDL_026"
$ws.Range("I27").Value = "select * from information_schema.schemata where schema_name in ('META','DINGO','ROOT','MYSQL','INFORMATION_SCHEMA','MYDDL_026')"
$ws.Range("J27").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_026.csv"
$ws.Range("K27").Value = "csv_containsAll"
$ws.Rows.Item(27).RowHeight = 27
$ws.Range("H27").WrapText = $true

# Row 28 (tall row, wrapped text)
$ws.Range("A28").Value = "ddl_027"
$ws.Range("B28").Value = "y"
$ws.Range("C28").Value = "删除schema,验证元数据信息删除"
$ws.Range("D28").Value = "DDL"
$ws.Range("E28").Value = "schemaDrop"
$ws.Range("H28").Value = "create schema MYDDL_027;drop schema MYDDL_027"
$ws.Range("I28").Value = "select * from information_schema.schemata where schema_name in ('META','DINGO','ROOT','MYSQL','INFORMATION_SCHEMA','MYDDL_027')"
$ws.Range("J28").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_027.csv"
$ws.Range("K28").Value = "csv_containsAll"
$ws.Rows.Item(28).RowHeight = 27
$ws.Range("H28").WrapText = $true

# Row 29
$ws.Range("A29").Value = "ddl_028"
$ws.Range("B29").Value = "y"
$ws.Range("C29").Value = "自定义创建的schema中创建的表查看元数据信息"
$ws.Range("D29").Value = "Schema"
$ws.Range("E29").Value = "Information_Schema"
$ws.Range("H29").Value = "create database MYDDL_028;create table MYDDL_028.MYDDL028_TBL01(id int, name varchar(20), primary key(id))"
$ws.Range("I29").Value = "select ``TABLE_CATALOG``,``TABLE_SCHEMA``,``TABLE_NAME``,``TABLE_TYPE``,``ENGINE``,``VERSION``,``ROW_FORMAT``,``TABLE_ROWS``,``AVG_ROW_LENGTH``,``DATA_LENGTH``,``MAX_DATA_LENGTH``,``INDEX_LENGTH``,``DATA_FREE``,``AUTO_INCREMENT``,``UPDATE_TIME``,``CHECK_TIME``,``TABLE_COLLATION``,``CHECKSUM``,``CREATE_OPTIONS``,``TABLE_COMMENT`` from information_schema.tables where ``TABLE_NAME``='MYDDL028_TBL01'"
$ws.Range("J29").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_028.csv"
$ws.Range("K29").Value = "csv_containsAll"

# Row 30
$ws.Range("A30").Value = "ddl_029"
$ws.Range("B30").Value = "y"
$ws.Range("C30").Value = "自定义创建的schema中创建的表删除后查看元数据信息"
$ws.Range("D30").Value = "Schema"
$ws.Range("E30").Value = "Information_Schema"
$ws.Range("H30").Value = "create database MYDDL_029;create table MYDDL_029.MYDDL029_TBL01(id int, name varchar(20), primary key(id));drop table MYDDL_029.MYDDL029_TBL01"
$ws.Range("I30").Value = "select ``TABLE_CATALOG``,``TABLE_SCHEMA``,``TABLE_NAME``,``TABLE_TYPE``,``ENGINE``,``VERSION``,``ROW_FORMAT``,``TABLE_ROWS``,``AVG_ROW_LENGTH``,``DATA_LENGTH``,``MAX_DATA_LENGTH``,``INDEX_LENGTH``,``DATA_FREE``,``AUTO_INCREMENT``,``UPDATE_TIME``,``CHECK_TIME``,``TABLE_COLLATION``,``CHECKSUM``,``CREATE_OPTIONS``,``TABLE_COMMENT`` from information_schema.tables where ``TABLE_NAME``='MYDDL029_TBL01'"
$ws.Range("J30").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_029.csv"
$ws.Range("K30").Value = "csv_containsAll"

# Row 31
$ws.Range("A31").Value = "ddl_030"
$ws.Range("B31").Value = "y"
$ws.Range("C31").Value = "MYSQL schema中创建的表查看元数据信息"
$ws.Range("D31").Value = "Schema"
$ws.Range("E31").Value = "Information_Schema"
$ws.Range("H31").Value = "create table MYSQL.MYDDL030_TBL01(id int not null auto_increment, name varchar(20), primary key(id))"
$ws.Range("I31").Value = "select ``TABLE_CATALOG``,``TABLE_SCHEMA``,``TABLE_NAME``,``TABLE_TYPE``,``ENGINE``,``VERSION``,``ROW_FORMAT``,``TABLE_ROWS``,``AVG_ROW_LENGTH``,``DATA_LENGTH``,``MAX_DATA_LENGTH``,``INDEX_LENGTH``,``DATA_FREE``,``AUTO_INCREMENT``,``UPDATE_TIME``,``CHECK_TIME``,``TABLE_COLLATION``,``CHECKSUM``,``CREATE_OPTIONS``,``TABLE_COMMENT`` from information_schema.tables where ``TABLE_NAME``='MYDDL030_TBL01'"
$ws.Range("J31").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_030.csv"
$ws.Range("K31").Value = "csv_containsAll"

# Row 32
$ws.Range("A32").Value = "ddl_031"
$ws.Range("B32").Value = "y"
$ws.Range("C32").Value = "MYSQL schema中创建的表删除后查看元数据信息"
$ws.Range("D32").Value = "Schema"
$ws.Range("E32").Value = "Information_Schema"
$ws.Range("H32").Value = "create table MYSQL.MYDDL031_TBL01(id int not null auto_increment, name varchar(20), primary key(id));drop table MYSQL.MYDDL031_TBL01"
$ws.Range("I32").Value = "select ``TABLE_CATALOG``,``TABLE_SCHEMA``,``TABLE_NAME``,``TABLE_TYPE``,``ENGINE``,``VERSION``,``ROW_FORMAT``,``TABLE_ROWS``,``AVG_ROW_LENGTH``,``DATA_LENGTH``,``MAX_DATA_LENGTH``,``INDEX_LENGTH``,``DATA_FREE``,``AUTO_INCREMENT``,``UPDATE_TIME``,``CHECK_TIME``,``TABLE_COLLATION``,``CHECKSUM``,``CREATE_OPTIONS``,``TABLE_COMMENT`` from information_schema.tables where ``TABLE_NAME``='MYDDL031_TBL01'"
$ws.Range("J32").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/ddl/expectedresult/ddl_031.csv"
$ws.Range("K32").Value = "csv_containsAll"

# ---- View settings ----
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("G27").Select()
